$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 2356
$ws.Range("E2").Value = -99
$ws.Range("F2").Value = -99
$ws.Range("G2").Value = -146
$ws.Range("H2").Value = -107
$ws.Range("I2").Value = -108
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 2782
$ws.Range("L2").Value = 566
$ws.Range("M2").Value = 2217
$ws.Range("N2").Value = 2221
$ws.Range("O2").Value = -4
$ws.Range("P2").Value = 290
$ws.Range("Q2").Value = 233
$ws.Range("R2").Value = -575
$ws.Range("S2").Value = 182
$ws.Range("T2").Value = 109
$ws.Range("U2").Value = 125
$ws.Range("V2").Value = 0
$ws.Range("W2").Value = -4.21
$ws.Range("X2").Value = -4.52
$ws.Range("Y2").Value = -5
$ws.Range("Z2").Value = -3.89
$ws.Range("AA2").Value = 25.52
$ws.Range("AB2").Value = 797.63
$ws.Range("AC2").Value = -186
$ws.Range("AD2").Value = -8.75
$ws.Range("AE2").Value = 3872
$ws.Range("AF2").Value = 0.42
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 57943763

# Row 3
$ws.Range("D3").Value = 2062
$ws.Range("E3").Value = -115
$ws.Range("F3").Value = -115
$ws.Range("G3").Value = -50
$ws.Range("H3").Value = -51
$ws.Range("I3").Value = -52
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 2403
$ws.Range("L3").Value = 371
$ws.Range("M3").Value = 2032
$ws.Range("N3").Value = 2036
$ws.Range("O3").Value = -4
$ws.Range("P3").Value = 290
$ws.Range("Q3").Value = 160
$ws.Range("R3").Value = 53
$ws.Range("S3").Value = -130
$ws.Range("T3").Value = 42
$ws.Range("U3").Value = 118
$ws.Range("V3").Value = 2
$ws.Range("W3").Value = -5.58
$ws.Range("X3").Value = -2.5
$ws.Range("Y3").Value = -2.44
$ws.Range("Z3").Value = -1.99
$ws.Range("AA3").Value = 18.24
$ws.Range("AB3").Value = 655.36
$ws.Range("AC3").Value = -90
$ws.Range("AD3").Value = -26.01
$ws.Range("AE3").Value = 3888
$ws.Range("AF3").Value = 0.6
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 57943763

# Row 4
$ws.Range("D4").Value = 1709
$ws.Range("E4").Value = -14
$ws.Range("F4").Value = -14
$ws.Range("G4").Value = 20
$ws.Range("H4").Value = 20
$ws.Range("I4").Value = 20
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 2423
$ws.Range("L4").Value = 390
$ws.Range("M4").Value = 2033
$ws.Range("N4").Value = 2037
$ws.Range("O4").Value = -4
$ws.Range("P4").Value = 290
$ws.Range("Q4").Value = 124
$ws.Range("R4").Value = -120
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 28
$ws.Range("U4").Value = 96
$ws.Range("V4").Value = 2
$ws.Range("W4").Value = -0.85
$ws.Range("X4").Value = 1.18
$ws.Range("Y4").Value = 0.98
$ws.Range("Z4").Value = 0.84
$ws.Range("AA4").Value = 19.17
$ws.Range("AB4").Value = 657.83
$ws.Range("AC4").Value = 34
$ws.Range("AD4").Value = 67.77
$ws.Range("AE4").Value = 3889
$ws.Range("AF4").Value = 0.6
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 57943763

# Row 5 (O5 gets fully cleared -- removed from the sheet)
$ws.Range("D5").Value = 1774
$ws.Range("E5").Value = 83
$ws.Range("F5").Value = 83
$ws.Range("G5").Value = 29
$ws.Range("H5").Value = 29
$ws.Range("I5").Value = 29
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 2382
$ws.Range("L5").Value = 321
$ws.Range("M5").Value = 2061
$ws.Range("N5").Value = 2061
$ws.Range("O5").ClearContents()
$ws.Range("P5").Value = 290
$ws.Range("Q5").Value = 47
$ws.Range("R5").Value = -79
$ws.Range("S5").Value = 0
$ws.Range("T5").Value = 99
$ws.Range("U5").Value = -52
$ws.Range("V5").Value = 2
$ws.Range("W5").Value = 4.67
$ws.Range("X5").Value = 1.65
$ws.Range("Y5").Value = 1.42
$ws.Range("Z5").Value = 1.22
$ws.Range("AA5").Value = 15.56
$ws.Range("AB5").Value = 668.42
$ws.Range("AC5").Value = 50
$ws.Range("AD5").Value = 55.46
$ws.Range("AE5").Value = 3937
$ws.Range("AF5").Value = 0.7
$ws.Range("AG5").Value = 50
$ws.Range("AH5").Value = 1.8
$ws.Range("AI5").Value = 90.31
$ws.Range("AJ5").Value = 57943763

# Row 6 (J6 and O6 remain absent, not touched)
$ws.Range("D6").Value = 1769
$ws.Range("E6").Value = 25
$ws.Range("F6").Value = 25
$ws.Range("G6").Value = 32
$ws.Range("H6").Value = 63
$ws.Range("I6").Value = 66
$ws.Range("K6").Value = 2621
$ws.Range("L6").Value = 377
$ws.Range("M6").Value = 2244
$ws.Range("N6").Value = 2225
$ws.Range("P6").Value = 290
$ws.Range("Q6").Value = 126
$ws.Range("R6").Value = -115
$ws.Range("S6").Value = 4
$ws.Range("T6").Value = 169
$ws.Range("U6").Value = -43
$ws.Range("V6").Value = 5
$ws.Range("W6").Value = 1.43
$ws.Range("X6").Value = 3.58
$ws.Range("Y6").Value = 3.07
$ws.Range("Z6").Value = 2.53
$ws.Range("AA6").Value = 16.81
$ws.Range("AB6").Value = 684.35
$ws.Range("AC6").Value = 113
$ws.Range("AD6").Value = 17.85
$ws.Range("AE6").Value = 4248
$ws.Range("AF6").Value = 0.48
$ws.Range("AG6").Value = 50
$ws.Range("AH6").Value = 2.47
$ws.Range("AI6").Value = 39.83
$ws.Range("AJ6").Value = 57943763

# Rows 7,8,9: clear all data cells (D..AJ), keep only A,B,C
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()
